$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 08:52"

# --- Ucrania (row 39) gets a full data refresh ---
$ws.Range("B39").Value = 10406
$ws.Range("C39").Value = 540
$ws.Range("D39").Value = 1238
$ws.Range("E39").Value = 8907
$ws.Range("F39").Value = 138
$ws.Range("G39").Value = 11
$ws.Range("H39").Value = 261

# --- Chequia (row 45) data refresh ---
$ws.Range("B45").Value = 7581
$ws.Range("C45").Value = 2
$ws.Range("D45").Value = 3120
$ws.Range("E45").Value = 4234

# --- Irak / Uzbekistan swap ranking (Uzbekistan overtakes Irak) ---
# Row 68 now holds Uzbekistan's refreshed figures, row 69 keeps Irak's
# (previously row 68's) figures unchanged.
$ws.Range("A68").Value = "Uzbekistan"
$ws.Range("B68").Value = 2017
$ws.Range("C68").Value = 15
$ws.Range("D68").Value = 1096
$ws.Range("E68").Value = 912
$ws.Range("F68").Value = 8
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 9

$ws.Range("A69").Value = "Irak"
$ws.Range("B69").Value = 2003
$ws.Range("C69").Value = 0
$ws.Range("D69").Value = 1346
$ws.Range("E69").Value = 565
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 92

# --- Letonia (row 93) data refresh ---
$ws.Range("B93").Value = 858
$ws.Range("C93").Value = 9
$ws.Range("E93").Value = 495
$ws.Range("F93").Value = 3

# --- Georgia (row 109) data refresh ---
$ws.Range("B109").Value = 539
$ws.Range("C109").Value = 22
$ws.Range("E109").Value = 355

# --- Islas Feroe (row 135) data refresh ---
$ws.Range("D135").Value = 184
$ws.Range("E135").Value = 3
